# Update "想去人数" (F) and "最低票价" (G) figures on the 展览 and 全部类型
# sheets to reflect the latest scrape output.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("F2").Value = 426
    $ws.Range("F3").Value = 5264
    $ws.Range("G3").Value = 62

    if ($name -eq "展览") {
        $ws.Range("F4").Value = 52
        $ws.Range("F5").Value = 57
        $ws.Range("F7").Value = 507
    }
    elseif ($name -eq "全部类型") {
        $ws.Range("F5").Value = 52
        $ws.Range("F6").Value = 57
        $ws.Range("F9").Value = 507
    }
}
